$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Hour 14 / Working on Assignment 2
$ws.Range("A14").Value = "Hour 14"
$ws.Range("B14").Value = "Working on Assignment 2"

# Update the selected cell to match the new active cell in the diff
$ws.Range("B14").Select()
